# Updated cryptos list data (Price and Volume(1h) columns)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range('D2')
$cell.NumberFormat = '@'
$cell.Value = '26.877.16'
$cell = $ws.Range('E2')
$cell.NumberFormat = '@'
$cell.Value = '  +0.23%  '
$cell = $ws.Range('D3')
$cell.NumberFormat = '@'
$cell.Value = '1.639.45'
$cell = $ws.Range('E3')
$cell.NumberFormat = '@'
$cell.Value = '  -0.20%  '
$cell = $ws.Range('E4')
$cell.NumberFormat = '@'
$cell.Value = '  -0.52%  '
$cell = $ws.Range('D5')
$cell.NumberFormat = '@'
$cell.Value = '216.87'
$cell = $ws.Range('E5')
$cell.NumberFormat = '@'
$cell.Value = '  -0.67%  '
$cell = $ws.Range('E6')
$cell.NumberFormat = '@'
$cell.Value = '  +1.80%  '
$cell = $ws.Range('E7')
$cell.NumberFormat = '@'
$cell.Value = '  -0.49%  '
$cell = $ws.Range('E8')
$cell.NumberFormat = '@'
$cell.Value = '  +1.60%  '
$cell = $ws.Range('E9')
$cell.NumberFormat = '@'
$cell.Value = '  +0.39%  '
$cell = $ws.Range('D10')
$cell.NumberFormat = '@'
$cell.Value = '19.87'
$cell = $ws.Range('E10')
$cell.NumberFormat = '@'
$cell.Value = '  +3.40%  '
$cell = $ws.Range('E11')
$cell.NumberFormat = '@'
$cell.Value = '  -0.08%  '
$cell = $ws.Range('D12')
$cell.NumberFormat = '@'
$cell.Value = '1.868.49'
$cell = $ws.Range('D13')
$cell.NumberFormat = '@'
$cell.Value = '1.645.14'
$cell = $ws.Range('E13')
$cell.NumberFormat = '@'
$cell.Value = '  +0.20%  '
$cell = $ws.Range('E14')
$cell.NumberFormat = '@'
$cell.Value = '  -0.77%  '
$cell = $ws.Range('D16')
$cell.NumberFormat = '@'
$cell.Value = '67.21'
$cell = $ws.Range('E16')
$cell.NumberFormat = '@'
$cell.Value = '  +3.15%  '
$cell = $ws.Range('D17')
$cell.NumberFormat = '@'
$cell.Value = '26.879.00'
$cell = $ws.Range('E17')
$cell.NumberFormat = '@'
$cell.Value = '  +0.24%  '
$cell = $ws.Range('E18')
$cell.NumberFormat = '@'
$cell.Value = '  -0.50%  '
$cell = $ws.Range('D19')
$cell.NumberFormat = '@'
$cell.Value = '219.51'
$cell = $ws.Range('E19')
$cell.NumberFormat = '@'
$cell.Value = '  +1.75%  '
$cell = $ws.Range('E20')
$cell.NumberFormat = '@'
$cell.Value = '  -0.51%  '
$cell = $ws.Range('D21')
$cell.NumberFormat = '@'
$cell.Value = '6.85'
$cell = $ws.Range('E21')
$cell.NumberFormat = '@'
$cell.Value = '  +3.13%  '
$cell = $ws.Range('E22')
$cell.NumberFormat = '@'
$cell.Value = '  +0.53%  '
$cell = $ws.Range('D23')
$cell.NumberFormat = '@'
$cell.Value = '2.44'
$cell = $ws.Range('E23')
$cell.NumberFormat = '@'
$cell.Value = '  +3.76%  '
$cell = $ws.Range('E24')
$cell.NumberFormat = '@'
$cell.Value = '  -0.37%  '
$cell = $ws.Range('D25')
$cell.NumberFormat = '@'
$cell.Value = '146.95'
$cell = $ws.Range('E25')
$cell.NumberFormat = '@'
$cell.Value = '  -0.45%  '
$cell = $ws.Range('E26')
$cell.NumberFormat = '@'
$cell.Value = '  -0.52%  '
$cell = $ws.Range('D27')
$cell.NumberFormat = '@'
$cell.Value = '7.35'
$cell = $ws.Range('E27')
$cell.NumberFormat = '@'
$cell.Value = '  +3.16%  '
$cell = $ws.Range('E28')
$cell.NumberFormat = '@'
$cell.Value = '  +0.44%  '
$cell = $ws.Range('D29')
$cell.NumberFormat = '@'
$cell.Value = '15.79'
$cell = $ws.Range('E29')
$cell.NumberFormat = '@'
$cell.Value = '  +0.31%  '
$cell = $ws.Range('E30')
$cell.NumberFormat = '@'
$cell.Value = '  -0.82%  '
$cell = $ws.Range('E31')
$cell.NumberFormat = '@'
$cell.Value = '  -0.89%  '
$cell = $ws.Range('E32')
$cell.NumberFormat = '@'
$cell.Value = '  -1.45%  '
$cell = $ws.Range('E33')
$cell.NumberFormat = '@'
$cell.Value = '  +0.69%  '
$cell = $ws.Range('E34')
$cell.NumberFormat = '@'
$cell.Value = '  +1.21%  '
$cell = $ws.Range('D35')
$cell.NumberFormat = '@'
$cell.Value = '1.269.56'
$cell = $ws.Range('E35')
$cell.NumberFormat = '@'
$cell.Value = '  +0.12%  '
$cell = $ws.Range('E36')
$cell.NumberFormat = '@'
$cell.Value = '  -0.31%  '
$cell = $ws.Range('E37')
$cell.NumberFormat = '@'
$cell.Value = '  +1.91%  '
$cell = $ws.Range('E38')
$cell.NumberFormat = '@'
$cell.Value = '  +0.27%  '
$cell = $ws.Range('E39')
$cell.NumberFormat = '@'
$cell.Value = '  +2.09%  '
$cell = $ws.Range('E40')
$cell.NumberFormat = '@'
$cell.Value = '  -0.51%  '
$cell = $ws.Range('E41')
$cell.NumberFormat = '@'
$cell.Value = '  +0.92%  '
$cell = $ws.Range('D42')
$cell.NumberFormat = '@'
$cell.Value = '5.40'
$cell = $ws.Range('E42')
$cell.NumberFormat = '@'
$cell.Value = '  +0.95%  '
$cell = $ws.Range('D43')
$cell.NumberFormat = '@'
$cell.Value = '1.779.33'
$cell = $ws.Range('E43')
$cell.NumberFormat = '@'
$cell.Value = '  -0.16%  '
$cell = $ws.Range('D44')
$cell.NumberFormat = '@'
$cell.Value = '2.10'
$cell = $ws.Range('E44')
$cell.NumberFormat = '@'
$cell.Value = '  -1.46%  '
$cell = $ws.Range('D45')
$cell.NumberFormat = '@'
$cell.Value = '61.79'
$cell = $ws.Range('E45')
$cell.NumberFormat = '@'
$cell.Value = '  +0.83%  '
$cell = $ws.Range('D46')
$cell.NumberFormat = '@'
$cell.Value = '91.84'
$cell = $ws.Range('E46')
$cell.NumberFormat = '@'
$cell.Value = '  -0.99%  '
$cell = $ws.Range('E47')
$cell.NumberFormat = '@'
$cell.Value = '  -0.88%  '
$cell = $ws.Range('E48')
$cell.NumberFormat = '@'
$cell.Value = '  +4.15%  '
$cell = $ws.Range('D49')
$cell.NumberFormat = '@'
$cell.Value = '0.0513'
$cell = $ws.Range('E49')
$cell.NumberFormat = '@'
$cell.Value = '  -0.48%  '
$cell = $ws.Range('D50')
$cell.NumberFormat = '@'
$cell.Value = '7.65'
$cell = $ws.Range('E50')
$cell.NumberFormat = '@'
$cell.Value = '  +1.18%  '
$cell = $ws.Range('E51')
$cell.NumberFormat = '@'
$cell.Value = '  -0.24%  '
